$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 274.91666
$ws.Range("J58").Value = 566.3333
$ws.Range("L58").Value = 1698.9999
$ws.Range("N58").Value = -1998.9999
$ws.Range("H68").Value = 174583.17
$ws.Range("J68").Value = 174583.17
$ws.Range("L68").Value = 174583.17
$ws.Range("N68").Value = -176081.17
$ws.Range("H70").Value = 1869.8
$ws.Range("I70").Value = 825
$ws.Range("J70").Value = 2566.3333
$ws.Range("K70").Value = 2475
$ws.Range("L70").Value = 7698.999899999999
$ws.Range("M70").Value = -2205
$ws.Range("N70").Value = -8238.999899999999
$ws.Range("H71").Value = 174583.17
$ws.Range("J71").Value = 174583.17
$ws.Range("L71").Value = 523749.51
$ws.Range("N71").Value = -531237.51
$ws.Range("H73").Value = 1869.8
$ws.Range("I73").Value = 825
$ws.Range("J73").Value = 2566.3333
$ws.Range("K73").Value = 2475
$ws.Range("L73").Value = 7698.999899999999
$ws.Range("M73").Value = -1539
$ws.Range("N73").Value = -9570.999899999999
$ws.Range("H86").Value = 4730.0312
$ws.Range("I86").Value = 2799.5
$ws.Range("K86").Value = 2799.5
$ws.Range("M86").Value = -1676.5
$ws.Range("H89").Value = 4730.0312
$ws.Range("I89").Value = 2799.5
$ws.Range("K89").Value = 13997.5
$ws.Range("M89").Value = -8381.5
$ws.Range("H92").Value = 567.43475
$ws.Range("I92").Value = 649.6667
$ws.Range("K92").Value = 649.6667
$ws.Range("M92").Value = 598.3333
$ws.Range("H106").Value = 22001894
$ws.Range("I106").Value = 25883892
$ws.Range("K106").Value = 25883892
$ws.Range("M106").Value = -25883261
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2522.5
$ws.Range("I45").Value = 1409.125
$ws.Range("J45").Value = 4749.25
$ws.Range("K45").Value = 1409.125
$ws.Range("L45").Value = 4749.25
$ws.Range("M45").Value = -1032.125
$ws.Range("N45").Value = -5503.25
$ws.Range("H74").Value = 3628.853
$ws.Range("I74").Value = 2008.7037
$ws.Range("J74").Value = 9878
$ws.Range("K74").Value = 2008.7037
$ws.Range("L74").Value = 9878
$ws.Range("M74").Value = -1134.7037
$ws.Range("N74").Value = -11626
$ws.Range("H77").Value = 3628.853
$ws.Range("I77").Value = 2008.7037
$ws.Range("J77").Value = 9878
$ws.Range("K77").Value = 10043.5185
$ws.Range("L77").Value = 49390
$ws.Range("M77").Value = -5675.5185
$ws.Range("N77").Value = -58126

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 776.6
$ws.Range("I64").Value = 907
$ws.Range("J64").Value = 689.6667
$ws.Range("K64").Value = 907
$ws.Range("L64").Value = 689.6667
$ws.Range("M64").Value = -682
$ws.Range("N64").Value = -1139.6667
$ws.Range("H67").Value = 776.6
$ws.Range("I67").Value = 907
$ws.Range("J67").Value = 689.6667
$ws.Range("K67").Value = 907
$ws.Range("L67").Value = 689.6667
$ws.Range("M67").Value = -127
$ws.Range("N67").Value = -2249.6667

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3029.4
$ws.Range("I3").Value = 3029.4
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 9088.200000000001
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -8976.200000000001
$ws.Range("N3").ClearContents()
$ws.Range("H23").Value = 95
$ws.Range("J23").Value = 102.333336
$ws.Range("L23").Value = 307.000008
$ws.Range("N23").Value = -777.000008
$ws.Range("H33").Value = 890.3333
$ws.Range("I33").Value = 646.6
$ws.Range("J33").Value = 1195
$ws.Range("K33").Value = 3879.6
$ws.Range("L33").Value = 7170
$ws.Range("M33").Value = -3596.6
$ws.Range("N33").Value = -7736
$ws.Range("H97").Value = 862
$ws.Range("J97").Value = 286.5
$ws.Range("L97").Value = 859.5
$ws.Range("N97").Value = -1851.5
$ws.Range("H114").Value = 840.5
$ws.Range("I114").Value = 456.5
$ws.Range("J114").Value = 1992.5
$ws.Range("K114").Value = 1369.5
$ws.Range("L114").Value = 5977.5
$ws.Range("M114").Value = 1884.5
$ws.Range("N114").Value = -12485.5
$ws.Range("H117").Value = 5999.5
$ws.Range("I117").Value = 1000
$ws.Range("J117").Value = 10999
$ws.Range("K117").Value = 3000
$ws.Range("L117").Value = 32997
$ws.Range("M117").Value = 442
$ws.Range("N117").Value = -39881

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1289.8889
$ws.Range("I2").Value = 1289.8889
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1289.8889
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -1176.8889
$ws.Range("N2").ClearContents()
$ws.Range("H80").Value = 9997.385
$ws.Range("I80").Value = 6685.1113
$ws.Range("K80").Value = 6685.1113
$ws.Range("M80").Value = -5687.1113
$ws.Range("H83").Value = 9997.385
$ws.Range("I83").Value = 6685.1113
$ws.Range("K83").Value = 33425.5565
$ws.Range("M83").Value = -28433.5565

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2299.182
$ws.Range("I22").Value = 1873.75
$ws.Range("J22").Value = 2542.2856
$ws.Range("K22").Value = 1873.75
$ws.Range("L22").Value = 2542.2856
$ws.Range("M22").Value = -1578.75
$ws.Range("N22").Value = -3132.2856
$ws.Range("H27").Value = 2299.182
$ws.Range("I27").Value = 1873.75
$ws.Range("J27").Value = 2542.2856
$ws.Range("K27").Value = 1873.75
$ws.Range("L27").Value = 2542.2856
$ws.Range("M27").Value = -1766.75
$ws.Range("N27").Value = -2756.2856
$ws.Range("H68").Value = 7913.5
$ws.Range("J68").Value = 4975
$ws.Range("L68").Value = 4975
$ws.Range("N68").Value = -6473
$ws.Range("H71").Value = 7913.5
$ws.Range("J71").Value = 4975
$ws.Range("L71").Value = 24875
$ws.Range("N71").Value = -32363
$ws.Range("H100").Value = 6009.8
$ws.Range("I100").Value = 6288.6665
$ws.Range("K100").Value = 6288.6665
$ws.Range("M100").Value = -5747.6665
$ws.Range("H132").Value = 8646.137000000001
$ws.Range("J132").Value = 11736.667
$ws.Range("L132").Value = 35210.001
$ws.Range("N132").Value = -40270.001
